$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue "D2" "27.772.81"
Set-TextValue "E2" "  +0.41%  "
Set-TextValue "D3" "1.848.37"
Set-TextValue "E3" "  -0.21%  "
Set-TextValue "D4" "1.006"
Set-TextValue "E4" "  -0.21%  "
Set-TextValue "D5" "334.83"
Set-TextValue "E5" "  +0.24%  "
Set-TextValue "E6" "  -0.24%  "
Set-TextValue "D7" "0.4645"
Set-TextValue "E7" "  +0.71%  "
Set-TextValue "D8" "0.3865"
Set-TextValue "E8" "  -0.75%  "
Set-TextValue "D9" "46.63"
Set-TextValue "E9" "  +1.10%  "
Set-TextValue "D10" "0.07900"
Set-TextValue "E10" "  -0.76%  "
Set-TextValue "D11" "0.9683"
Set-TextValue "E11" "  -3.94%  "
Set-TextValue "D12" "21.29"
Set-TextValue "E12" "  -1.61%  "
Set-TextValue "B13" "Polkadot"
Set-TextValue "C13" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D13" "5.886"
Set-TextValue "E13" "  -1.66%  "
Set-TextValue "B14" "WrappedEther"
Set-TextValue "C14" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D14" "1.827.41"
Set-TextValue "E14" "  -2.17%  "
Set-TextValue "D15" "7.156"
Set-TextValue "E15" "  -0.50%  "
Set-TextValue "D16" "1.006"
Set-TextValue "E16" "  -0.23%  "
Set-TextValue "D17" "90.06"
Set-TextValue "E17" "  +1.73%  "
Set-TextValue "D18" "0.06605"
Set-TextValue "E18" "  -1.43%  "
Set-TextValue "E19" "  -0.95%  "
Set-TextValue "D20" "17.35"
Set-TextValue "E20" "  +0.57%  "
Set-TextValue "D21" "1.004"
Set-TextValue "E21" "  -0.50%  "
Set-TextValue "D22" "27.757.57"
Set-TextValue "E22" "  +0.42%  "
Set-TextValue "D23" "5.342"
Set-TextValue "E23" "  -1.38%  "
Set-TextValue "E24" "  -1.45%  "
Set-TextValue "D25" "2.293"
Set-TextValue "E25" "  -0.85%  "
Set-TextValue "D26" "158.57"
Set-TextValue "E26" "  -0.50%  "
Set-TextValue "D27" "2.048.19"
Set-TextValue "E27" "  -1.65%  "
Set-TextValue "D28" "19.49"
Set-TextValue "E28" "  -0.53%  "
Set-TextValue "D29" "2.068"
Set-TextValue "E29" "  -2.87%  "
Set-TextValue "D30" "5.359"
Set-TextValue "E30" "  -1.85%  "
Set-TextValue "D31" "118.74"
Set-TextValue "E31" "  -1.98%  "
Set-TextValue "D32" "0.09412"
Set-TextValue "E32" "  -0.05%  "
Set-TextValue "D33" "0.9450"
Set-TextValue "E33" "  -3.56%  "
Set-TextValue "D34" "3.585"
Set-TextValue "E34" "  -0.94%  "
Set-TextValue "D35" "5.263"
Set-TextValue "E35" "  -0.95%  "
Set-TextValue "D36" "1.328"
Set-TextValue "E36" "  -1.51%  "
Set-TextValue "D37" "0.06003"
Set-TextValue "E37" "  -0.42%  "
Set-TextValue "D38" "0.02209"
Set-TextValue "E38" "  -1.09%  "
Set-TextValue "D39" "8.241"
Set-TextValue "E39" "  -1.12%  "
Set-TextValue "E40" "  -0.24%  "
Set-TextValue "D41" "1.159"
Set-TextValue "E41" "  -2.27%  "
Set-TextValue "D42" "0.5807"
Set-TextValue "E42" "  -2.16%  "
Set-TextValue "D43" "0.1840"
Set-TextValue "E43" "  -1.61%  "
Set-TextValue "D44" "10.11"
Set-TextValue "E44" "  -2.70%  "
Set-TextValue "D45" "1.279"
Set-TextValue "E45" "  +2.89%  "
Set-TextValue "B46" "EnergySwap"
Set-TextValue "C46" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D46" "11.99"
Set-TextValue "E46" "  -1.54%  "
Set-TextValue "B47" "Decentraland"
Set-TextValue "C47" "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue "D47" "0.5445"
Set-TextValue "E47" "  -2.74%  "
Set-TextValue "D48" "1.928"
Set-TextValue "E48" "  +0.52%  "
Set-TextValue "D49" "0.06839"
Set-TextValue "E49" "  +1.83%  "
Set-TextValue "D50" "110.90"
Set-TextValue "E50" "  -0.27%  "
Set-TextValue "E51" "  -32.79%  "
